$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new columns D:G for rows 2-5 (sex, environ, sire, dam)
$ws.Range("D2").Value = "U"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

$ws.Range("D3").Value = "U"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2

$ws.Range("D4").Value = "U"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2

$ws.Range("D5").Value = "U"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 3

# Update the active selection to match the saved workbook state
$ws.Range("E9").Select()
